$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 54, shifting existing rows 54-147 down to 57-150
$ws.Rows("54:56").Insert()

# Populate the 3 newly inserted rows with their data
# Row 54
$ws.Range("A54").Value = 8
$ws.Range("B54").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C54").Value = 'Coquimbo'
$ws.Range("D54").Value = 44544
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112021
$ws.Range("G54").Value = 'Ají'
$ws.Range("H54").Value = 'Inferno'
$ws.Range("I54").Value = 'Primera'
$ws.Range("J54").Value = 500
$ws.Range("K54").Value = 13000
$ws.Range("L54").Value = 14000
$ws.Range("M54").Value = 13500
$ws.Range("N54").Value = '$/caja 12 kilos'
$ws.Range("O54").Value = 'Región de Arica y Parinacota'
$ws.Range("P54").Value = 1125
$ws.Range("Q54").Value = 12
$ws.Range("R54").Value = 'Hortaliza'

# Row 55
$ws.Range("A55").Value = 8
$ws.Range("B55").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C55").Value = 'Coquimbo'
$ws.Range("D55").Value = 44544
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = 100112021
$ws.Range("G55").Value = 'Ají'
$ws.Range("H55").Value = 'Inferno'
$ws.Range("I55").Value = 'Primera'
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 17000
$ws.Range("L55").Value = 18000
$ws.Range("M55").Value = 17500
$ws.Range("N55").Value = '$/caja 15 kilos'
$ws.Range("O55").Value = 'Provincia de Limarí'
$ws.Range("P55").Value = 1167
$ws.Range("Q55").Value = 15
$ws.Range("R55").Value = 'Hortaliza'

# Row 56
$ws.Range("A56").Value = 8
$ws.Range("B56").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C56").Value = 'Coquimbo'
$ws.Range("D56").Value = 44544
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 100112021
$ws.Range("G56").Value = 'Ají'
$ws.Range("H56").Value = 'Inferno'
$ws.Range("I56").Value = 'Segunda'
$ws.Range("J56").Value = 360
$ws.Range("K56").Value = 8000
$ws.Range("L56").Value = 8500
$ws.Range("M56").Value = 8250
$ws.Range("N56").Value = '$/caja 12 kilos'
$ws.Range("O56").Value = 'Región de Arica y Parinacota'
$ws.Range("P56").Value = 688
$ws.Range("Q56").Value = 12
$ws.Range("R56").Value = 'Hortaliza'

